$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("C1").Value = "Assessment Type"
$ws.Range("D1").Value = "Marks "
$ws.Range("E1").Value = "Maximim Marks"
$ws.Range("F1").Value = "Date "
$ws.Range("G1").Value = "Semester "

# Add new data row (row 2)
$ws.Range("A2").Value = "STU576"
$ws.Range("B2").Value = "Ssc"
$ws.Range("C2").Value = "class assesment"
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 200
$ws.Range("F2").Value = "'07/08/2025"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "5th"
